$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# Grab references to all the shapes we need to touch *first*, while
# their original (unique) names are still intact -- a couple of the
# renames below collide with each other's old/new names, so resolving
# every lookup by name up-front avoids picking the wrong shape later.
$oldPic = Get-ShapeByName $s "Picture 1"
$tb1 = Get-ShapeByName $s "TextBox 21"
$tb2 = Get-ShapeByName $s "TextBox 22"
$tb3 = Get-ShapeByName $s "TextBox 9"

# ---------------------------------------------------------------------
# 1. Duplicate the existing background picture, move the duplicate to
#    the back of the z-order (right after grpSpPr) with a nudged
#    position, then delete the original picture. This reproduces the
#    diff's "new <p:pic> inserted at top (rId2) + old <p:pic> (rId4)
#    removed at bottom" behaviour, since Shape.Duplicate() re-embeds the
#    picture under a brand-new relationship id.
# ---------------------------------------------------------------------
$newPic = $oldPic.Duplicate()
$newPic.ZOrder(1)
$newPic.Left = -2.88
$newPic.Top = 5.184566929133858
$newPic.Width = 645.1200259999999
$newPic.Height = 483.84001199999994
$oldPic.Delete()

# ---------------------------------------------------------------------
# 2. (a) label -> "TextBox 9", reposition/resize, retext as
#    "(a) Passive Inference" at 22pt.
# ---------------------------------------------------------------------
$tb1.Name = "TextBox 9"
$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "(a"
$tr1.Font.Size = 22
$run1b = $tr1.InsertAfter(") ")
$run1b.Font.Size = 22
$run1c = $tr1.InsertAfter("Passive Inference")
$run1c.Font.Size = 22
$tb1.Left = 71.88055118110236
$tb1.Top = 219.79866841732283
$tb1.Width = 234.39701087401573
$tb1.Height = 33.92811023622047

# ---------------------------------------------------------------------
# 3. (b) label -> "TextBox 10", reposition/resize, retext as
#    "(b) Active Inference" at 22pt.
# ---------------------------------------------------------------------
$tb2.Name = "TextBox 10"
$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "(b"
$tr2.Font.Size = 22
$run2b = $tr2.InsertAfter(") ")
$run2b.Font.Size = 22
$run2c = $tr2.InsertAfter("Active Inference")
$run2c.Font.Size = 22
$tb2.Left = 378.15811023622047
$tb2.Top = 219.72212998425195
$tb2.Width = 231.56346956692911
$tb2.Height = 33.92811023622047

# ---------------------------------------------------------------------
# 4. (c) label -> "TextBox 11", reposition/resize, retext as
#    "(c) Beliefs in true goal" at 22pt.
# ---------------------------------------------------------------------
$tb3.Name = "TextBox 11"
$tr3 = $tb3.TextFrame.TextRange
$tr3.Text = "(c"
$tr3.Font.Size = 22
$run3b = $tr3.InsertAfter(") Beliefs in true goal")
$run3b.Font.Size = 22
$tb3.Left = 200.63503937007874
$tb3.Top = 465.7434645669291
$tb3.Width = 243.84992225984251
$tb3.Height = 33.92811023622047
